$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new user row (row 3) — fixes the "loading" bug where a newly
# registered user wasn't being written to the sheet.
$ws.Range("A3").Value = "test1"
$ws.Range("B3").Value = "test1234"
$ws.Range("C3").Value = "test@gmail.com"
$ws.Range("D3").Value = "test"
$ws.Range("E3").Value = "test"
$ws.Range("F3").Value = "test"
$ws.Range("G3").Value = 123456

$ws.Range("A3:G3").Select()
